$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.267.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.029.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.75%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.61%  "

$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  +4.41%  "

$ws.Range("E10").Value = "  +2.99%  "

$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.323.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.850"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.83%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.028.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.253.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0862"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.89%  "

$ws.Range("E21").Value = "  +3.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.46%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("E24").Value = "  +5.27%  "

$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.42%  "

$ws.Range("E28").Value = "  -4.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.97%  "

$ws.Range("E30").Value = "  +5.89%  "

$ws.Range("E31").Value = "  +1.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0671"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.12%  "

$ws.Range("E34").Value = "  +14.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.52"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.60%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.21%  "

$ws.Range("E37").Value = "  -0.16%  "

$ws.Range("E38").Value = "  +1.45%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0971"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0217"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.37%  "

$ws.Range("E43").Value = "  +1.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.377.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("E48").Value = "  +2.97%  "

$ws.Range("E49").Value = "  +15.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.37%  "
